$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update stats for 2025-08 (row 21)
$ws.Range("B21").Value = 6255
$ws.Range("D21").Value = 5658845
$ws.Range("E21").Value = 904.691446842526
$ws.Range("F21").Value = 8.574900190939072
$ws.Range("H21").Value = 29.16870959230051
